# Re-apply the latest scraped cryptos.xlsx snapshot onto the sheet.
# (GitHub Actions refresh: prices/volumes updated, a couple of rows swapped.)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking price strings that must stay TEXT (they use
# '.' as both thousands- and decimal-separator, so Excel would silently coerce
# them to numbers and drop formatting like trailing zeros). A leading apostrophe
# is the standard Excel "force text" entry marker and is stripped from the
# stored value, so the cell ends up holding the exact digit string below.

$ws.Range('D2').Value = '29.877.72'
$ws.Range('E2').Value = '  +0.13%  '
$ws.Range('D3').Value = '1.894.23'
$ws.Range('E3').Value = '  +0.04%  '
$ws.Range('D4').Value = '''1.001'
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').Value = '''0.7736'
$ws.Range('E5').Value = '  -1.95%  '
$ws.Range('D6').Value = '''244.24'
$ws.Range('E6').Value = '  +0.41%  '
$ws.Range('D7').Value = '''1.000'
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('E8').Value = '  -0.46%  '
$ws.Range('D9').Value = '''25.67'
$ws.Range('E9').Value = '  +1.60%  '
$ws.Range('D10').Value = '''0.07226'
$ws.Range('E10').Value = '  -0.20%  '
$ws.Range('D11').Value = '''0.08864'
$ws.Range('E11').Value = '  +9.53%  '
$ws.Range('B12').Value = 'Polygon'
$ws.Range('C12').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D12').Value = '''0.7711'
$ws.Range('E12').Value = '  +0.93%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.950.89'
$ws.Range('E13').Value = '  +2.61%  '
$ws.Range('D14').Value = '''5.435'
$ws.Range('E14').Value = '  -1.83%  '
$ws.Range('E15').Value = '  +2.28%  '
$ws.Range('D16').Value = '''6.196'
$ws.Range('E16').Value = '  +0.91%  '
$ws.Range('D17').Value = '29.995.00'
$ws.Range('E17').Value = '  +0.45%  '
$ws.Range('D18').Value = '''13.94'
$ws.Range('E18').Value = '  +0.45%  '
$ws.Range('D19').Value = '''245.49'
$ws.Range('E19').Value = '  +0.70%  '
$ws.Range('D20').Value = '''0.000007867'
$ws.Range('E20').Value = '  +1.22%  '
$ws.Range('D21').Value = '2.194.65'
$ws.Range('E21').Value = '  +1.25%  '
$ws.Range('D22').Value = '''8.160'
$ws.Range('E22').Value = '  +0.63%  '
$ws.Range('D23').Value = '''1.0000'
$ws.Range('E23').Value = '  -0.12%  '
$ws.Range('E24').Value = '  -0.16%  '
$ws.Range('D25').Value = '''0.1601'
$ws.Range('E25').Value = '  -2.77%  '
$ws.Range('D26').Value = '''9.512'
$ws.Range('E26').Value = '  +1.48%  '
$ws.Range('D27').Value = '''162.43'
$ws.Range('E27').Value = '  -0.62%  '
$ws.Range('D28').Value = '''18.80'
$ws.Range('E28').Value = '  +0.60%  '
$ws.Range('D29').Value = '''2.042'
$ws.Range('E29').Value = '  -0.22%  '
$ws.Range('E30').Value = '  +2.07%  '
$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D31').Value = '''1.542'
$ws.Range('E31').Value = '  -0.24%  '
$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D32').Value = '''4.552'
$ws.Range('E32').Value = '  +1.98%  '
$ws.Range('D33').Value = '''4.114'
$ws.Range('E33').Value = '  +0.60%  '
$ws.Range('D34').Value = '''0.05490'
$ws.Range('E34').Value = '  -0.83%  '
$ws.Range('D35').Value = '''1.248'
$ws.Range('E35').Value = '  -1.36%  '
$ws.Range('D36').Value = '''0.7519'
$ws.Range('E36').Value = '  +1.94%  '
$ws.Range('D37').Value = '''0.9998'
$ws.Range('E37').Value = '  +0.05%  '
$ws.Range('D38').Value = '''2.716'
$ws.Range('E38').Value = '  +3.52%  '
$ws.Range('D39').Value = '''0.01958'
$ws.Range('E39').Value = '  +1.82%  '
$ws.Range('D40').Value = '''2.789'
$ws.Range('E40').Value = '  +0.31%  '
$ws.Range('D41').Value = '''0.4507'
$ws.Range('E41').Value = '  +2.00%  '
$ws.Range('D42').Value = '''73.84'
$ws.Range('E42').Value = '  -0.06%  '
$ws.Range('D43').Value = '1.091.62'
$ws.Range('E43').Value = '  -4.67%  '
$ws.Range('D44').Value = '''6.032'
$ws.Range('E44').Value = '  +2.75%  '
$ws.Range('D45').Value = '''0.8553'
$ws.Range('E45').Value = '  +0.71%  '
$ws.Range('E46').Value = '  -0.05%  '
$ws.Range('E47').Value = '  +0.75%  '
$ws.Range('D48').Value = '''102.77'
$ws.Range('E48').Value = '  -1.36%  '
$ws.Range('D49').Value = '''7.613'
$ws.Range('E49').Value = '  +2.40%  '
$ws.Range('D50').Value = '''9.850'
$ws.Range('E50').Value = '  -1.87%  '
$ws.Range('D51').Value = '2.067.83'
$ws.Range('E51').Value = '  +0.44%  '
